# feat: add 2022-Q1 data
#
# Before: sheet tabs = [ "2021-Q4", "总计" ]
# After:  sheet tabs = [ "2021-Q4", "2022-Q1", "总计" ]
#
# The existing "总计" (grand-total) sheet is duplicated in place so the
# original keeps its summary-table formatting/shared-string pool. The
# original is renamed to "2022-Q1" and repurposed to hold the per-fund
# holding detail for the new quarter; the duplicate is renamed back to
# "总计" and gets a new leading row for "2022-Q1" (pushing "2021-Q4" down).

$wb = $excel.ActiveWorkbook

$zongji = $wb.Worksheets.Item("总计")

# Duplicate "总计" right after itself -- the copy inherits identical
# formatting/styles, which is what both resulting sheets need.
$zongji.Copy($null, $zongji)

$q1    = $wb.Worksheets.Item(2)   # was "总计", becomes "2022-Q1"
$total = $wb.Worksheets.Item(3)   # the fresh copy, becomes "总计"

$q1.Name    = "2022-Q1"
$total.Name = "总计"

# ---------------------------------------------------------------------
# "2022-Q1": replace the 3-column summary layout with the 7-column
# per-fund holding-detail layout (same shape as the "2021-Q4" sheet).
# ---------------------------------------------------------------------

# New header cells E1:H1 need the same bold/bordered look as B1:H1 --
# grab it from the existing B1 header cell before overwriting any text.
$q1.Range("B1").Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$q1.Range("A2").Value = 0

$q1.Range("B2:G2").NumberFormat = "@"
$q1.Range("B2").Value = "512040"
$q1.Range("C2").Value = "富国中证价值ETF"
$q1.Range("D2").Value = "3.44"
$q1.Range("E2").Value = "99.55"
$q1.Range("F2").Value = "1.40"
$q1.Range("G2").Value = "0.0482"
$q1.Range("B2:G2").Style = "Normal"

$q1.Range("H2").Value = 4

# ---------------------------------------------------------------------
# "总计": push the existing "2021-Q4" row down to row 3 and write the
# new "2022-Q1" row into row 2 (row 3 already carries over correct
# data/format from the duplicated sheet, so only its index cell A3
# and row 2 need attention).
# ---------------------------------------------------------------------

$total.Range("A2").Copy()
$total.Range("A3").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.05

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 1
$total.Range("D3").Value = 0
